# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) values and a refreshed cover
# image URL to the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---------------------------------------------------------
$wsExhibit.Range("F2").Value  = 198
$wsExhibit.Range("F3").Value  = 5458
$wsExhibit.Range("I4").Value  = "//i2.hdslb.com/bfs/openplatform/202404/ukD6OzH51713324745614.jpeg"
$wsExhibit.Range("F7").Value  = 635
$wsExhibit.Range("F8").Value  = 606
$wsExhibit.Range("F12").Value = 4754
$wsExhibit.Range("F13").Value = 447
$wsExhibit.Range("F14").Value = 211
$wsExhibit.Range("F17").Value = 3618
$wsExhibit.Range("F18").Value = 190
$wsExhibit.Range("F19").Value = 1126
$wsExhibit.Range("F21").Value = 44
$wsExhibit.Range("F26").Value = 145
$wsExhibit.Range("F28").Value = 329
$wsExhibit.Range("F31").Value = 23
$wsExhibit.Range("F32").Value = 36
$wsExhibit.Range("F33").Value = 37

# --- Sheet "全部类型" ------------------------------------------------------
$wsAll.Range("F2").Value  = 198
$wsAll.Range("F4").Value  = 5458
$wsAll.Range("I5").Value  = "//i2.hdslb.com/bfs/openplatform/202404/ukD6OzH51713324745614.jpeg"
$wsAll.Range("F8").Value  = 635
$wsAll.Range("F9").Value  = 606
$wsAll.Range("F13").Value = 4754
$wsAll.Range("F14").Value = 447
$wsAll.Range("F15").Value = 211
$wsAll.Range("F18").Value = 3619
$wsAll.Range("F19").Value = 190
$wsAll.Range("F20").Value = 1126
$wsAll.Range("F22").Value = 44
$wsAll.Range("F27").Value = 145
$wsAll.Range("F29").Value = 329
$wsAll.Range("F32").Value = 23
$wsAll.Range("F33").Value = 36
$wsAll.Range("F34").Value = 37

$wb.Save()
